$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9: update title
$ws.Range("D9").Value = "[공지] SIAI 설립 관련 서류 공개 및 음해성 공격 해명"

# Row 37: update title and link
$ws.Range("D37").Value = "[Paper Review] MixText: Linguistically-Informed Interpolation of Hidden Space for Semi-Supervised Text Classification"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1925&mod=document&pageid=1"

# Row 41: update title and link
$ws.Range("D41").Value = "상향식 기술 구현과 현실의 모방"
$ws.Range("E41").Value = "http://cloudinsight.net/cloud/%ec%83%81%ed%96%a5%ec%8b%9d-%ea%b8%b0%ec%88%a0-%ea%b5%ac%ed%98%84%ea%b3%bc-%ed%98%84%ec%8b%a4%ec%9d%98-%eb%aa%a8%eb%b0%a9/"

# Row 52: update title and link
$ws.Range("D52").Value = "텍스트 데이터 파일 읽기에 계속 실패할 경우"
$ws.Range("E52").Value = "http://ds.sumeun.org/?p=2514&utm_source=rss&utm_medium=rss&utm_campaign=%25ed%2585%258d%25ec%258a%25a4%25ed%258a%25b8-%25eb%258d%25b0%25ec%259d%25b4%25ed%2584%25b0-%25ed%258c%258c%25ec%259d%25bc-%25ec%259d%25bd%25ea%25b8%25b0%25ec%2597%2590-%25ea%25b3%2584%25ec%2586%258d-%25ec%258b%25a4%25ed%258c%25a8%25ed%2595%25a0-%25ea%25b2%25bd%25ec%259a%25b0"
